$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = '26.484.29'
$ws.Cells.Item(2,4).Style = "Normal"

# Row 3
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = '1.677.81'
$ws.Cells.Item(3,5).Value = '  +3.44%  '
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Style = "Normal"

# Row 4
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value = '  +0.05%  '
$ws.Cells.Item(4,5).Style = "Normal"

# Row 5
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '216.47'
$ws.Cells.Item(5,4).Style = "Normal"

# Row 6
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '0.5317'
$ws.Cells.Item(6,5).Value = '  +2.19%  '
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Style = "Normal"

# Row 7
$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value = '  +0.02%  '
$ws.Cells.Item(7,5).Style = "Normal"

# Row 8
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = '  +4.32%  '
$ws.Cells.Item(8,5).Style = "Normal"

# Row 9
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '0.06393'
$ws.Cells.Item(9,5).Value = '  +2.07%  '
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Style = "Normal"

# Row 10
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '21.60'
$ws.Cells.Item(10,5).Value = '  +5.85%  '
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Style = "Normal"

# Row 11
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '0.07801'
$ws.Cells.Item(11,5).Value = '  +3.54%  '
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Style = "Normal"

# Row 12
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,5).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '1.682.21'
$ws.Cells.Item(12,5).Value = '  +3.52%  '
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Style = "Normal"

# Row 13
$ws.Cells.Item(13,5).NumberFormat = "@"
$ws.Cells.Item(13,5).Value = '  +2.78%  '
$ws.Cells.Item(13,5).Style = "Normal"

# Row 14
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '0.5563'
$ws.Cells.Item(14,5).Value = '  +1.91%  '
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Style = "Normal"

# Row 15
$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value = '  +5.34%  '
$ws.Cells.Item(15,5).Style = "Normal"

# Row 16
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,5).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '65.64'
$ws.Cells.Item(16,5).Value = '  +2.22%  '
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Style = "Normal"

# Row 17
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,5).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '26.518.73'
$ws.Cells.Item(17,5).Value = '  +2.62%  '
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Style = "Normal"

# Row 18
$ws.Cells.Item(18,5).NumberFormat = "@"
$ws.Cells.Item(18,5).Value = '  +0.00%  '
$ws.Cells.Item(18,5).Style = "Normal"

# Row 19
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,5).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '4.759'
$ws.Cells.Item(19,5).Value = '  +2.37%  '
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Style = "Normal"

# Row 20
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,5).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '194.65'
$ws.Cells.Item(20,5).Value = '  +5.83%  '
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Style = "Normal"

# Row 21
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,5).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '10.36'
$ws.Cells.Item(21,5).Value = '  +3.01%  '
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Style = "Normal"

# Row 22
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '6.337'
$ws.Cells.Item(22,5).Value = '  +4.47%  '
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Style = "Normal"

# Row 23
$ws.Cells.Item(23,5).NumberFormat = "@"
$ws.Cells.Item(23,5).Value = '  +0.05%  '
$ws.Cells.Item(23,5).Style = "Normal"

# Row 24
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,5).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '143.85'
$ws.Cells.Item(24,5).Value = '  -0.65%  '
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Style = "Normal"

# Row 25
$ws.Cells.Item(25,5).NumberFormat = "@"
$ws.Cells.Item(25,5).Value = '  +5.94%  '
$ws.Cells.Item(25,5).Style = "Normal"

# Row 26
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,5).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '7.432'
$ws.Cells.Item(26,5).Value = '  +1.02%  '
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Style = "Normal"

# Row 27
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,5).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '16.32'
$ws.Cells.Item(27,5).Value = '  +4.95%  '
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Style = "Normal"

# Row 28
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,5).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '1.427'
$ws.Cells.Item(28,5).Value = '  +5.16%  '
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Style = "Normal"

# Row 29
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,5).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '0.06155'
$ws.Cells.Item(29,5).Value = '  +4.87%  '
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Style = "Normal"

# Row 30
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,5).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '1.274'
$ws.Cells.Item(30,5).Value = '  +2.90%  '
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Style = "Normal"

# Row 31
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,5).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '3.618'
$ws.Cells.Item(31,5).Value = '  +6.70%  '
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Style = "Normal"

# Row 32
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,5).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '3.454'
$ws.Cells.Item(32,5).Value = '  +3.12%  '
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Style = "Normal"

# Row 33
$ws.Cells.Item(33,5).NumberFormat = "@"
$ws.Cells.Item(33,5).Value = '  +5.16%  '
$ws.Cells.Item(33,5).Style = "Normal"

# Row 34
$ws.Cells.Item(34,5).NumberFormat = "@"
$ws.Cells.Item(34,5).Value = '  +3.46%  '
$ws.Cells.Item(34,5).Style = "Normal"

# Row 35
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,5).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '2.430'
$ws.Cells.Item(35,5).Value = '  +2.00%  '
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Style = "Normal"

# Row 36
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,5).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = '2.785'
$ws.Cells.Item(36,5).Value = '  +2.37%  '
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Style = "Normal"

# Row 37
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,5).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '0.5742'
$ws.Cells.Item(37,5).Value = '  -0.27%  '
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Style = "Normal"

# Row 38
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,5).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '0.01637'
$ws.Cells.Item(38,5).Value = '  +3.06%  '
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Style = "Normal"

# Row 39
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,5).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = '6.033'
$ws.Cells.Item(39,5).Value = '  +6.63%  '
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Style = "Normal"

# Row 40
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,5).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '1.072.74'
$ws.Cells.Item(40,5).Value = '  +4.58%  '
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Style = "Normal"

# Row 41
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,5).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.8598'
$ws.Cells.Item(41,5).Value = '  +1.78%  '
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Style = "Normal"

# Row 42
$ws.Cells.Item(42,5).NumberFormat = "@"
$ws.Cells.Item(42,5).Value = '  -0.24%  '
$ws.Cells.Item(42,5).Style = "Normal"

# Row 43
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,5).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '100.03'
$ws.Cells.Item(43,5).Value = '  +0.64%  '
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Style = "Normal"

# Row 44
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,5).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '1.824.87'
$ws.Cells.Item(44,5).Value = '  +3.10%  '
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Style = "Normal"

# Row 45
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,5).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '57.01'
$ws.Cells.Item(45,5).Value = '  +4.47%  '
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Style = "Normal"

# Row 46
$ws.Cells.Item(46,2).Value = 'EnergySwap'
$ws.Cells.Item(46,3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,5).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '8.138'
$ws.Cells.Item(46,5).Value = '  +2.67%  '
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Style = "Normal"

# Row 47
$ws.Cells.Item(47,2).Value = 'Frax'
$ws.Cells.Item(47,3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,5).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '1.002'
$ws.Cells.Item(47,5).Value = '  +0.41%  '
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Style = "Normal"

# Row 48
$ws.Cells.Item(48,2).Value = 'BabyDogeCoin'
$ws.Cells.Item(48,3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,5).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '0.0₈104'
$ws.Cells.Item(48,5).Value = '  -5.09%  '
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Style = "Normal"

# Row 49
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,5).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '0.05207'
$ws.Cells.Item(49,5).Value = '  +1.06%  '
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Style = "Normal"

# Row 50
$ws.Cells.Item(50,2).Value = 'Aptos'
$ws.Cells.Item(50,3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,5).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '6.033'
$ws.Cells.Item(50,5).Value = '  +3.73%  '
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Style = "Normal"

# Row 51
$ws.Cells.Item(51,2).Value = 'Mantle'
$ws.Cells.Item(51,3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,5).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '0.4237'
$ws.Cells.Item(51,5).Value = '  +0.68%  '
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Style = "Normal"
